# Apply the "Update SALN Reports Generation" change:
# Add a new "AS OF" date column (column J) to the PERM-CAS-COT sheet,
# recording the date each filer's SALN was received/processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PERM-CAS-COT")

# --- Column width for the new column J ---
$ws.Columns.Item(10).ColumnWidth = 16.25

# --- Header cells (row 4/5), copy formatting from the existing I4:I5 header ---
$ws.Range("I4:I5").Copy()
$ws.Range("J4:J5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J4").Value = "AS OF"
$ws.Range("J4:J5").Merge()

# --- Data cells: copy formatting from column I on each data row, then set date + format ---
$dataRows = @(7, 8, 9, 10, 12, 13, 15, 16, 17, 18, 19, 20)
$dates = @{
    7  = 45308
    8  = 45308
    9  = 45308
    10 = 45303
    12 = 45308
    13 = 45307
    15 = 45308
    16 = 45308
    17 = 45308
    18 = 45310
    19 = 45310
    20 = 45291
}

foreach ($r in $dataRows) {
    $ws.Range("I$r").Copy()
    $ws.Range("J$r").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("J$r").Value = $dates[$r]
    $ws.Range("J$r").NumberFormat = "mmm d, yyyy"
}

$excel.CutCopyMode = 0
